# Append three new paragraphs (TECHNICAL SKILLS heading + two "Languages" lines)
# after the existing content, matching the target OOXML exactly.  We build the
# insertion as a WordprocessingML fragment and use Range.InsertXML so the new
# runs get exactly the rPr/pPr we want, with nothing inherited from the
# paragraph mark currently at the end of the document.

$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)

$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="2" w:color="black"/></w:pBdr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000" w:themeColor="accent1"/><w:sz w:val="24"/></w:rPr><w:t>TECHNICAL SKILLS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000" w:themeColor="accent1"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">   Languages : </w:t></w:r><w:r><w:t>Python,C,JavaScript...etc</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">   Languages : </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000" w:themeColor="accent1"/><w:sz w:val="24"/></w:rPr><w:t>Python,C,JavaScript...etc</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$endRange.InsertXML($xmlPayload)
Write-Host "Inserted TECHNICAL SKILLS + Languages paragraphs; paragraph count now:" $d.Paragraphs.Count
